# Apply text replacements to the division-problems worksheet.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-12-24 Tuesday" "2024-12-25 Wednesday"

Replace-Text "50÷3=" "27÷8="
Replace-Text "48÷5=" "78÷5="
Replace-Text "87÷8=" "55÷5="
Replace-Text "31÷2=" "84÷2="
Replace-Text "69÷4=" "47÷5="

Replace-Text "62÷9=" "51÷8="
Replace-Text "45÷8=" "67÷9="
Replace-Text "33÷5=" "28÷7="
Replace-Text "25÷5=" "24÷2="
Replace-Text "22÷7=" "62÷2="

Replace-Text "37÷5=" "22÷6="
Replace-Text "51÷2=" "50÷7="
Replace-Text "49÷5=" "77÷7="
Replace-Text "84÷8=" "41÷6="
Replace-Text "84÷9=" "13÷9="

Replace-Text "38÷5=" "43÷6="
Replace-Text "66÷2=" "32÷4="
Replace-Text "30÷6=" "36÷9="
Replace-Text "11÷2=" "24÷4="
Replace-Text "80÷4=" "83÷8="

Replace-Text "20÷9=" "95÷7="
Replace-Text "21÷6=" "19÷8="
Replace-Text "34÷4=" "34÷9="
Replace-Text "29÷7=" "46÷4="
Replace-Text "79÷8=" "68÷4="
